# Adding labels/technologyLabel to output Excel file in createFileString()
#
# Relabels the "parameter" column entries to the underscore-joined
# identifiers used by createFileString(), updates the "Updated" date,
# and refreshes a handful of assumption values on the "Inputs of Model"
# sheet. Finally moves the active selection to B12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Updated" date (B3) -------------------------------------------------
# Typing a dd.mm.yyyy-looking string directly gets auto-parsed into a date
# serial by the smart-input path, which would also swap the cell's style.
# Going through a formula and then Paste-Special (values only) keeps the
# original style (s="5") while still landing a literal text value.
$ws.Range("B3").Formula = "=""06.10.2020"""
$ws.Range("B3").Copy()
[void]$ws.Range("B3").PasteSpecial(-4163)

# --- parameter labels (column A) -----------------------------------------
$ws.Range("A7").Value  = "Minimum_daily_mileage"
$ws.Range("A8").Value  = "Battery_capacity"
$ws.Range("A9").Value  = "Electric_consumption_NEFZ"
$ws.Range("A10").Value = "Fuel_consumption_NEFZ"
$ws.Range("A11").Value = "Electric_consumption_Artemis"
$ws.Range("A12").Value = "Fuel_consumption_Artemis"
$ws.Range("A13").Value = "Maximum_SOC"
$ws.Range("A14").Value = "Minimum_SOC"
$ws.Range("A15").Value = "Rated_power_of_charging_column"
$ws.Range("A16").Value = "Is_BEV?"

# --- updated assumption values (column B) ---------------------------------
$ws.Range("B8").Value  = 17
$ws.Range("B9").Value  = 20
$ws.Range("B11").Value = 19
$ws.Range("B15").Value = 11

# --- selection -------------------------------------------------------------
[void]$ws.Range("B12").Select()
